$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "http://10.0.74.5/index/home"
$ws.Range("B2").Value = "4fku01"
$ws.Range("D2").Value = "17/10/2019"
$ws.Range("E2").Value = "13:39:52.016"
$ws.Range("F2").Value = "17/10/2019"
$ws.Range("G2").Value = "13:40:45.020"
$ws.Range("I2").Value = "Firefox 69.0.3"
